# This script applies the weekly refresh of "Fruta / hortaliza" price data
# for the "Hortaliza, Femacal de La Calera - Perejil" sheet:
#  - the Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
#    Precio promedio ponderado (M) and Precio $/Kg (P) values for most
#    existing data rows are refreshed with new figures
#  - one new data row (31) is appended at the bottom of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with refreshed Fecha / Volumen / Precio values
$ws.Range("D2").Value = 44223 ; $ws.Range("J2").Value = 80 ; $ws.Range("K2").Value = 2500 ; $ws.Range("L2").Value = 3000 ; $ws.Range("M2").Value = 2781 ; $ws.Range("P2").Value = 927
$ws.Range("D3").Value = 44627 ; $ws.Range("J3").Value = 78 ; $ws.Range("K3").Value = 3500 ; $ws.Range("L3").Value = 3500 ; $ws.Range("M3").Value = 3500 ; $ws.Range("P3").Value = 1167
$ws.Range("D4").Value = 44225 ; $ws.Range("J4").Value = 56 ; $ws.Range("K4").Value = 3000 ; $ws.Range("M4").Value = 3000 ; $ws.Range("P4").Value = 1000
$ws.Range("D5").Value = 44557 ; $ws.Range("J5").Value = 104 ; $ws.Range("K5").Value = 2000 ; $ws.Range("L5").Value = 2500 ; $ws.Range("M5").Value = 2260 ; $ws.Range("P5").Value = 753
$ws.Range("D6").Value = 44340 ; $ws.Range("J6").Value = 54 ; $ws.Range("K6").Value = 3000 ; $ws.Range("L6").Value = 3000 ; $ws.Range("M6").Value = 3000 ; $ws.Range("P6").Value = 1000
$ws.Range("D7").Value = 44537 ; $ws.Range("J7").Value = 88 ; $ws.Range("L7").Value = 2200 ; $ws.Range("M7").Value = 2091 ; $ws.Range("P7").Value = 697
$ws.Range("D8").Value = 44224 ; $ws.Range("J8").Value = 67 ; $ws.Range("K8").Value = 3000 ; $ws.Range("M8").Value = 3000 ; $ws.Range("P8").Value = 1000
$ws.Range("D9").Value = 44165 ; $ws.Range("J9").Value = 68
$ws.Range("D10").Value = 44845 ; $ws.Range("J10").Value = 80 ; $ws.Range("K10").Value = 2500 ; $ws.Range("L10").Value = 2500 ; $ws.Range("M10").Value = 2500 ; $ws.Range("P10").Value = 833
$ws.Range("D12").Value = 44222 ; $ws.Range("J12").Value = 45
$ws.Range("D13").Value = 44804 ; $ws.Range("J13").Value = 85
$ws.Range("D14").Value = 44193 ; $ws.Range("J14").Value = 70
$ws.Range("D15").Value = 44574 ; $ws.Range("J15").Value = 50
$ws.Range("D16").Value = 44389 ; $ws.Range("J16").Value = 81 ; $ws.Range("K16").Value = 2800 ; $ws.Range("M16").Value = 2889 ; $ws.Range("P16").Value = 963
$ws.Range("D17").Value = 44260 ; $ws.Range("J17").Value = 60 ; $ws.Range("K17").Value = 3500 ; $ws.Range("L17").Value = 3500 ; $ws.Range("M17").Value = 3500 ; $ws.Range("P17").Value = 1167
$ws.Range("D18").Value = 44242 ; $ws.Range("J18").Value = 95 ; $ws.Range("K18").Value = 2500 ; $ws.Range("M18").Value = 2737 ; $ws.Range("P18").Value = 912
$ws.Range("D19").Value = 44187 ; $ws.Range("J19").Value = 65
$ws.Range("D20").Value = 44166 ; $ws.Range("J20").Value = 45 ; $ws.Range("L20").Value = 2500 ; $ws.Range("M20").Value = 2500 ; $ws.Range("P20").Value = 833
$ws.Range("D21").Value = 44390 ; $ws.Range("J21").Value = 50 ; $ws.Range("K21").Value = 3000 ; $ws.Range("L21").Value = 3000 ; $ws.Range("M21").Value = 3000 ; $ws.Range("P21").Value = 1000
$ws.Range("D23").Value = 44937 ; $ws.Range("J23").Value = 68 ; $ws.Range("K23").Value = 3500 ; $ws.Range("L23").Value = 3500 ; $ws.Range("M23").Value = 3500 ; $ws.Range("P23").Value = 1167
$ws.Range("D24").Value = 44179 ; $ws.Range("J24").Value = 78 ; $ws.Range("K24").Value = 3000 ; $ws.Range("M24").Value = 3000 ; $ws.Range("P24").Value = 1000
$ws.Range("D26").Value = 44669 ; $ws.Range("J26").Value = 92 ; $ws.Range("K26").Value = 2500 ; $ws.Range("L26").Value = 3000 ; $ws.Range("M26").Value = 2755 ; $ws.Range("P26").Value = 918
$ws.Range("D27").Value = 44292 ; $ws.Range("J27").Value = 40 ; $ws.Range("K27").Value = 3000 ; $ws.Range("L27").Value = 3000 ; $ws.Range("M27").Value = 3000 ; $ws.Range("P27").Value = 1000
$ws.Range("D28").Value = 44559 ; $ws.Range("J28").Value = 68 ; $ws.Range("K28").Value = 2000 ; $ws.Range("L28").Value = 2000 ; $ws.Range("M28").Value = 2000 ; $ws.Range("P28").Value = 667
$ws.Range("D29").Value = 44536 ; $ws.Range("J29").Value = 125 ; $ws.Range("K29").Value = 2200 ; $ws.Range("L29").Value = 2200 ; $ws.Range("M29").Value = 2200 ; $ws.Range("P29").Value = 733
$ws.Range("D30").Value = 44756 ; $ws.Range("J30").Value = 104 ; $ws.Range("K30").Value = 2800 ; $ws.Range("M30").Value = 2904 ; $ws.Range("P30").Value = 968

# Append the new data row (31) at the bottom of the table
$ws.Range("A31").Value = 3
$ws.Range("B31").Value = "Femacal de La Calera"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44291
$ws.Range("D31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = 100112044
$ws.Range("G31").Value = "Perejil"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 45
$ws.Range("K31").Value = 3000
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = 3000
$ws.Range("N31").Value = "$/docena de atados (3 kilos)"
$ws.Range("O31").Value = "Provincia de Quillota"
$ws.Range("P31").Value = 1000
$ws.Range("Q31").Value = 3
$ws.Range("R31").Value = "Hortaliza"
